# Updated symbol list on Sat Feb  4 22:45:06 UTC 2023 with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) figures for the
# crypto ticker rows that moved since the last scrape.
# NumberFormat is set to "@" (Text) before writing each value so that
# numeric-looking strings (e.g. "331.63", "0.28%") are kept as literal
# text instead of being auto-converted to numbers/percentages by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range('D2').Value = '331.63'
$ws.Range("E2").NumberFormat = "@"
$ws.Range('E2').Value = '0.28%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range('D3').Value = '41.35'
$ws.Range("E3").NumberFormat = "@"
$ws.Range('E3').Value = '0.16%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range('D4').Value = '5.690'
$ws.Range("E4").NumberFormat = "@"
$ws.Range('E4').Value = '-0.46%'
$ws.Range("E5").NumberFormat = "@"
$ws.Range('E5').Value = '4.19%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '8.815'
$ws.Range("E6").NumberFormat = "@"
$ws.Range('E6').Value = '1.04%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range('D7').Value = '4.504'
$ws.Range("E7").NumberFormat = "@"
$ws.Range('E7').Value = '-0.17%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range('D8').Value = '1.980'
$ws.Range("E8").NumberFormat = "@"
$ws.Range('E8').Value = '-2.82%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '0.9284'
$ws.Range("E10").NumberFormat = "@"
$ws.Range('E10').Value = '0.60%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range('D11').Value = '0.1251'
$ws.Range("E11").NumberFormat = "@"
$ws.Range('E11').Value = '0.32%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range('D12').Value = '0.1983'
$ws.Range("E12").NumberFormat = "@"
$ws.Range('E12').Value = '1.99%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range('D13').Value = '0.09535'
$ws.Range("E13").NumberFormat = "@"
$ws.Range('E13').Value = '2.39%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range('D14').Value = '0.03959'
$ws.Range("E14").NumberFormat = "@"
$ws.Range('E14').Value = '8.05%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range('D15').Value = '0.1063'
$ws.Range("E15").NumberFormat = "@"
$ws.Range('E15').Value = '0.71%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range('D16').Value = '0.001316'
$ws.Range("E16").NumberFormat = "@"
$ws.Range('E16').Value = '0.82%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range('D17').Value = '0.006106'
$ws.Range("E17").NumberFormat = "@"
$ws.Range('E17').Value = '-0.51%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range('D18').Value = '3.435'
$ws.Range("E18").NumberFormat = "@"
$ws.Range('E18').Value = '1.59%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range('E19').Value = '0.84%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range('D20').Value = '9.164'
$ws.Range("E20").NumberFormat = "@"
$ws.Range('E20').Value = '9.93%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range('E21').Value = '-3.75%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range('D22').Value = '0.2509'
$ws.Range("E22").NumberFormat = "@"
$ws.Range('E22').Value = '-5.38%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '0.04404'
$ws.Range("E23").NumberFormat = "@"
$ws.Range('E23').Value = '-0.58%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range('D24').Value = '0.001245'
$ws.Range("E24").NumberFormat = "@"
$ws.Range('E24').Value = '-1.35%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range('D25').Value = '0.004376'
$ws.Range("E25").NumberFormat = "@"
$ws.Range('E25').Value = '0.63%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '0.0001190'
$ws.Range("E26").NumberFormat = "@"
$ws.Range('E26').Value = '-4.08%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range('D27').Value = '0.0003990'
$ws.Range("E27").NumberFormat = "@"
$ws.Range('E27').Value = '-0.06%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range('D39').Value = '0.02831'
$ws.Range("E39").NumberFormat = "@"
$ws.Range('E39').Value = '0.87%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range('D40').Value = '0.05539'
$ws.Range("E40").NumberFormat = "@"
$ws.Range('E40').Value = '1.14%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range('D41').Value = '0.007900'
$ws.Range("E41").NumberFormat = "@"
$ws.Range('E41').Value = '3.46%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '0.1441'
$ws.Range("E42").NumberFormat = "@"
$ws.Range('E42').Value = '1.34%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range('E43').Value = '-10.16%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range('D44').Value = '0.002081'
$ws.Range("E44").NumberFormat = "@"
$ws.Range('E44').Value = '-1.80%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '0.01029'
$ws.Range("E45").NumberFormat = "@"
$ws.Range('E45').Value = '-13.36%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '0.00007308'
$ws.Range("E46").NumberFormat = "@"
$ws.Range('E46').Value = '8.04%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '0.00000000750'
$ws.Range("E47").NumberFormat = "@"
$ws.Range('E47').Value = '-0.04%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '0.003209'
$ws.Range("E48").NumberFormat = "@"
$ws.Range('E48').Value = '4.47%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range('D49').Value = '0.002278'
$ws.Range("E49").NumberFormat = "@"
$ws.Range('E49').Value = '-0.09%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '0.00002101'
$ws.Range("E50").NumberFormat = "@"
$ws.Range('E50').Value = '-0.04%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range('D51').Value = '0.0002001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range('E51').Value = '-0.04%'
